# SnappCarr P2P Car Rental App - Issues workbook update
#
# 1. Delete the stray "Sheet1" worksheet (it only held two helper lookup
#    values that are no longer needed).
# 2. Fill in Story Points (C) / "time spent"-style (D) estimate values for
#    the first block of tasks under "Setup" (rows 5-16) on the "Tasks"
#    sheet.
# 3. Move the on-screen selection down to where the user was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

# --- 1. Remove the extra "Sheet1" worksheet -------------------------------
if ($wb.Worksheets.Count -gt 1) {
    $excel.DisplayAlerts = $false
    $wb.Worksheets.Item("Sheet1").Delete() | Out-Null
    $excel.DisplayAlerts = $true
}

# --- 2. Populate the Story Points (C) and Hours (D) columns ---------------
$estimates = @{
    5  = @(1, 1)
    6  = @(1, 1)
    7  = @(1, 0.5)
    8  = @(1, 0.5)
    9  = @(1, 0.5)
    10 = @(1, 0.5)
    11 = @(1, 0.5)
    12 = @(1, 0.5)
    13 = @(1, 0.5)
    14 = @(1, 0.5)
    15 = @(1, 0.5)
    16 = @(1, 0.5)
}

foreach ($row in $estimates.Keys) {
    $vals = $estimates[$row]
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
}

# --- 3. Update the active selection ----------------------------------------
$ws.Activate()
$ws.Range("E26").Select() | Out-Null
